$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# "Gestionar mis vídeos" (old row 25) is replaced by a longer "Mi espacio" sub-menu.
# Insert 8 new rows right after the current row 25 so the old rows 26-27
# (9.0/Ajustes, 10.0/Contacto) end up at rows 34-35.
$ws.Rows("26:33").Insert()

# Update B25 in place, then fill the newly inserted rows 26-33.
$ws.Range("B25").Value = "Inicio"

$ws.Range("A26").Value = "8.2"
$ws.Range("B26").Value = "Mis comentarios"

$ws.Range("A27").Value = "8.3"
$ws.Range("B27").Value = "Historial"

$ws.Range("A28").Value = "8.4"
$ws.Range("B28").Value = "Mis listas de reproducción"

$ws.Range("A29").Value = "8.5"
$ws.Range("B29").Value = "Ajustes de mi espacio"

$ws.Range("A30").Value = "8.6"
$ws.Range("B30").Value = "Vídeos favoritos"

$ws.Range("A31").Value = "8.7"
$ws.Range("B31").Value = "Vídeos subidos"

$ws.Range("A32").Value = "8.8"
$ws.Range("B32").Value = "Mis suscripciones"

$ws.Range("A33").Value = "8.9"
$ws.Range("B33").Value = "Subir vídeo"

# Match the widened column B and the view left scrolled down to the new rows.
$ws.Columns("B:B").ColumnWidth = 21.08984375
$ws.Application.ActiveWindow.ScrollRow = 19
$ws.Range("B31").Select()
